# Updated cryptos list on Mon Dec 18 20:23:01 UTC 2023 with GitHub Actions
#
# Applies the per-cell value updates (prices / 1h volume deltas), including
# the Monero/Toncoin row-content swap at rows 30-31, to Sheet1 of the
# cryptos workbook. Cell values are plain text (inlineStr) in the source
# file, so every write forces a text NumberFormat while it's in flight and
# restores the cell's original Style afterwards -- this keeps numeric-
# looking strings (e.g. "0.604", "41.703.87") stored as text instead of
# being auto-coerced to Number by Excel's usual type inference, without
# leaving any stray style/number-format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = 'D2'; Value = '41.703.87' }
    @{ Cell = 'E2'; Value = '  -1.23%  ' }
    @{ Cell = 'D3'; Value = '2.172.47' }
    @{ Cell = 'E3'; Value = '  -3.03%  ' }
    @{ Cell = 'E4'; Value = '  +0.13%  ' }
    @{ Cell = 'D5'; Value = '237.85' }
    @{ Cell = 'E5'; Value = '  -2.10%  ' }
    @{ Cell = 'D6'; Value = '0.604' }
    @{ Cell = 'E6'; Value = '  -3.62%  ' }
    @{ Cell = 'D7'; Value = '72.21' }
    @{ Cell = 'E7'; Value = '  -2.59%  ' }
    @{ Cell = 'E8'; Value = '  -0.10%  ' }
    @{ Cell = 'D9'; Value = '0.580' }
    @{ Cell = 'E9'; Value = '  -3.48%  ' }
    @{ Cell = 'D10'; Value = '39.84' }
    @{ Cell = 'E10'; Value = '  -5.74%  ' }
    @{ Cell = 'D11'; Value = '0.0906' }
    @{ Cell = 'E11'; Value = '  -5.24%  ' }
    @{ Cell = 'D12'; Value = '54.52' }
    @{ Cell = 'E12'; Value = '  -3.64%  ' }
    @{ Cell = 'E13'; Value = '  -2.61%  ' }
    @{ Cell = 'D14'; Value = '6.70' }
    @{ Cell = 'E14'; Value = '  -3.38%  ' }
    @{ Cell = 'D15'; Value = '2.504.23' }
    @{ Cell = 'E15'; Value = '  -2.79%  ' }
    @{ Cell = 'E16'; Value = '  -0.01%  ' }
    @{ Cell = 'D17'; Value = '2.168.29' }
    @{ Cell = 'E17'; Value = '  -3.67%  ' }
    @{ Cell = 'E18'; Value = '  -7.26%  ' }
    @{ Cell = 'D19'; Value = '41.575.28' }
    @{ Cell = 'E19'; Value = '  -1.22%  ' }
    @{ Cell = 'D20'; Value = '0.0000102' }
    @{ Cell = 'E20'; Value = '  -2.23%  ' }
    @{ Cell = 'D21'; Value = '70.02' }
    @{ Cell = 'E21'; Value = '  -3.72%  ' }
    @{ Cell = 'D22'; Value = '5.78' }
    @{ Cell = 'E22'; Value = '  -7.08%  ' }
    @{ Cell = 'D23'; Value = '9.98' }
    @{ Cell = 'E23'; Value = '  -11.98%  ' }
    @{ Cell = 'D24'; Value = '225.90' }
    @{ Cell = 'E24'; Value = '  -1.86%  ' }
    @{ Cell = 'E25'; Value = '  +0.03%  ' }
    @{ Cell = 'E26'; Value = '  -0.07%  ' }
    @{ Cell = 'D27'; Value = '10.68' }
    @{ Cell = 'E27'; Value = '  -6.34%  ' }
    @{ Cell = 'E28'; Value = '  -9.95%  ' }
    @{ Cell = 'E29'; Value = '  -3.78%  ' }
    @{ Cell = 'B30'; Value = 'Monero' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D30'; Value = '171.59' }
    @{ Cell = 'E30'; Value = '  +2.68%  ' }
    @{ Cell = 'B31'; Value = 'Toncoin' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = 'D31'; Value = '2.08' }
    @{ Cell = 'E31'; Value = '  -5.58%  ' }
    @{ Cell = 'D32'; Value = '19.85' }
    @{ Cell = 'E32'; Value = '  -3.68%  ' }
    @{ Cell = 'D33'; Value = '33.04' }
    @{ Cell = 'E33'; Value = '  +11.07%  ' }
    @{ Cell = 'D34'; Value = '0.0772' }
    @{ Cell = 'E34'; Value = '  -4.57%  ' }
    @{ Cell = 'D35'; Value = '5.27' }
    @{ Cell = 'E35'; Value = '  -6.70%  ' }
    @{ Cell = 'E36'; Value = '  -3.66%  ' }
    @{ Cell = 'D37'; Value = '4.28' }
    @{ Cell = 'E37'; Value = '  -1.26%  ' }
    @{ Cell = 'E38'; Value = '  -7.22%  ' }
    @{ Cell = 'E39'; Value = '  +1.11%  ' }
    @{ Cell = 'D40'; Value = '11.97' }
    @{ Cell = 'E40'; Value = '  -9.05%  ' }
    @{ Cell = 'E41'; Value = '  -2.21%  ' }
    @{ Cell = 'D42'; Value = '5.36' }
    @{ Cell = 'E42'; Value = '  -6.26%  ' }
    @{ Cell = 'D43'; Value = '58.90' }
    @{ Cell = 'E43'; Value = '  -8.76%  ' }
    @{ Cell = 'D44'; Value = '0.189' }
    @{ Cell = 'E44'; Value = '  -5.04%  ' }
    @{ Cell = 'D45'; Value = '8.41' }
    @{ Cell = 'E45'; Value = '  -3.62%  ' }
    @{ Cell = 'D46'; Value = '0.0964' }
    @{ Cell = 'E46'; Value = '  -3.99%  ' }
    @{ Cell = 'D47'; Value = '97.06' }
    @{ Cell = 'E47'; Value = '  -7.29%  ' }
    @{ Cell = 'E48'; Value = '  -5.36%  ' }
    @{ Cell = 'D49'; Value = '1.11' }
    @{ Cell = 'E49'; Value = '  -4.90%  ' }
    @{ Cell = 'D50'; Value = '2.19' }
    @{ Cell = 'E50'; Value = '  -5.96%  ' }
    @{ Cell = 'D51'; Value = '2.62' }
    @{ Cell = 'E51'; Value = '  -2.04%  ' }
)

foreach ($edit in $edits) {
    $rng = $ws.Range($edit.Cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $edit.Value
    $rng.Style = $origStyle
}
